# Gendata.xlsx - "Add files via upload" edit
# - Decrement the "from"/"to" bus numbers (columns B and C) by 1 for rows 2-11
#   on the "Lines" sheet.
# - Switch the active/selected sheet from "Gen slack" to "Lines".
# - Update the Lines sheet selection from K25 to E18.

$wb = $excel.ActiveWorkbook

$linesSheet = $wb.Worksheets.Item("Lines")

# Decrement bus-number columns B and C (rows 2-11) by 1.
for ($row = 2; $row -le 11; $row++) {
    $bCell = $linesSheet.Cells.Item($row, 2)
    $cCell = $linesSheet.Cells.Item($row, 3)
    $bCell.Value = $bCell.Value2 - 1
    $cCell.Value = $cCell.Value2 - 1
}

# Make "Lines" the active sheet (this also flips tabSelected off "Gen slack"
# and onto "Lines", and updates the workbook's activeTab).
$linesSheet.Activate()

# Move the selection on the Lines sheet to E18.
$null = $linesSheet.Range("E18").Select()
